$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the definition text in C18
$ws.Range("C18").Value = ""

# Update B27 from "cryosphere" to "land ice"
$ws.Range("B27").Value = "land ice"

# Add "cryosphere" to B28
$ws.Range("B28").Value = "cryosphere"

# Add "plant ecology" to B29
$ws.Range("B29").Value = "plant ecology"

# Add "ocean ecology" to B30
$ws.Range("B30").Value = "ocean ecology"

# New rows 31-39
$data = @(
    @{ row = 31; a = "nicest-2-subjects:10013"; b = "biogeochemistry"; e = $null },
    @{ row = 32; a = "nicest-2-subjects:10014"; b = "chemistry"; e = $null },
    @{ row = 33; a = "nicest-2-subjects:10015"; b = "geospatial domain"; e = $null },
    @{ row = 34; a = "nicest-2-subjects:10016"; b = "energy"; e = $null },
    @{ row = 35; a = "nicest-2-subjects:10017"; b = "bioprospecting"; e = $null },
    @{ row = 36; a = "nicest-2-subjects:10018"; b = "prediction"; e = $null },
    @{ row = 37; a = "nicest-2-subjects:10019"; b = "historical"; e = $null },
    @{ row = 38; a = "nicest-2-subjects:10020"; b = "aerosols"; e = $null },
    @{ row = 39; a = "nicest-2-subjects:10021"; b = "Earth System Modelling"; e = "ESM" }
)

foreach ($item in $data) {
    $r = $item.row
    $ws.Cells.Item($r, 1).Value = $item.a
    $ws.Cells.Item($r, 2).Value = $item.b
    if ($item.e -ne $null) {
        $ws.Cells.Item($r, 5).Value = $item.e
    }
}
